$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.542.33"
$ws.Range("E2").Value = "  -4.26%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.323.28"
$ws.Range("E3").Value = "  -4.95%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.22"
$ws.Range("E5").Value = "  -3.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.10"
$ws.Range("E6").Value = "  -3.36%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.320.30"
$ws.Range("E8").Value = "  -4.98%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.481"
$ws.Range("E9").Value = "  -0.73%  "

$ws.Range("E10").Value = "  -4.21%  "

$ws.Range("E11").Value = "  -3.90%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.379"
$ws.Range("E12").Value = "  -2.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.882.62"
$ws.Range("E13").Value = "  -5.08%  "

$ws.Range("E14").Value = "  -0.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.312.68"
$ws.Range("E15").Value = "  -5.22%  "

$ws.Range("E16").Value = "  -5.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.90"
$ws.Range("E17").Value = "  +1.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.562.66"
$ws.Range("E18").Value = "  -4.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.59"
$ws.Range("E19").Value = "  +0.61%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.69"
$ws.Range("E20").Value = "  -1.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.08"
$ws.Range("E21").Value = "  -9.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "355.01"
$ws.Range("E22").Value = "  -7.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.558"
$ws.Range("E23").Value = "  -3.57%  "

$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.450.22"
$ws.Range("E25").Value = "  -5.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.76"
$ws.Range("E26").Value = "  -6.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000109"
$ws.Range("E27").Value = "  -5.27%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").Value = "  -0.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.24"
$ws.Range("E29").Value = "  -0.59%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.46"
$ws.Range("E30").Value = "  -1.60%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.89"
$ws.Range("E31").Value = "  -2.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.12"
$ws.Range("E32").Value = "  -5.51%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.152"
$ws.Range("E33").Value = "  -2.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.351.04"
$ws.Range("E35").Value = "  -4.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.70"
$ws.Range("E36").Value = "  -2.33%  "

$ws.Range("E37").Value = "  -2.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.86"
$ws.Range("E38").Value = "  -0.36%  "

$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.48"
$ws.Range("E39").Value = "  -3.16%  "

$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "160.88"
$ws.Range("E40").Value = "  -2.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0764"
$ws.Range("E41").Value = "  -2.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  -0.16%  "

$ws.Range("E43").Value = "  -0.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.14"
$ws.Range("E44").Value = "  -1.79%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.744"
$ws.Range("E45").Value = "  -7.77%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.12"
$ws.Range("E46").Value = "  -5.13%  "

$ws.Range("E47").Value = "  -4.57%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.43"
$ws.Range("E48").Value = "  -7.14%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.75"
$ws.Range("E49").Value = "  -0.58%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.866"
$ws.Range("E50").Value = "  -5.80%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.42"
$ws.Range("E51").Value = "  +2.59%  "
